{"js": "// Apply each literal text replacement by searching for the exact\n// 'before' run text and replacing it in place, leaving sibling runs\n// (e.g. <w:br/> line breaks) untouched.\nconst replacements = [\n  { find: \"Ativa\u00e7\u00e3o: 01/01/2011\", replace: \"Ativa\u00e7\u00e3o: 01/01/2023\" },\n  { find: \"Curso (semestre ideal): EM (4)\", replace: \"Curso (semestre ideal): EF (4), EM (4)\" },\n  { find: \"O objetivo desta disciplina \u00e9 apresentar ao aluno fundamentos de F\u00edsico-Qu\u00edmica, mais especificamente Termodin\u00e2mica Qu\u00edmica, com foco em problemas encontrados na Engenharia de Materiais. Bastante \u00eanfase \u00e9 dada nos assuntos relativos a Termoqu\u00edmica e Equil\u00edbrio Qu\u00edmica, abordando equil\u00edbrio com fases condensadas.\", replace: \"Esta disciplina visa apresentar fundamentos de termodin\u00e2mica aplicada \u00e0 \u00e1rea de ci\u00eancia e engenharia de materiais. Especial \u00eanfase \u00e9 dada \u00e0 energia na forma e calor para aquecimento de sistemas termodin\u00e2micos; c\u00e1lculos de varia\u00e7\u00e3o de entalpia; entropia e energia de Gibbs de elementos e compostos em mudan\u00e7as de estado; c\u00e1lculos de varia\u00e7\u00e3o de entalpia; entropia e energia de Gibbs de rea\u00e7\u00e3o; aplica\u00e7\u00e3o da propriedade  energia de Gibbs para avalia\u00e7\u00e3o de transforma\u00e7\u00f5es espont\u00e2neas e em equil\u00edbrio; fundamentos de termodin\u00e2mica de solu\u00e7\u00f5es; c\u00e1lculos de condi\u00e7\u00f5es de equil\u00edbrio em sistemas heterog\u00eaneos. Apresenta-se tamb\u00e9m as principais diferen\u00e7as entre esta disciplina e a disciplina de Termodin\u00e2mica de M\u00e1quinas.\" },\n  { find: \"1) 1a Lei da Termodin\u00e2mica2) 2a e 3a Leis da Termodin\u00e2mica3) Rela\u00e7\u00f5es entre Propriedades Termodin\u00e2micas4) Equil\u00edbrio5) Equil\u00edbrio Qu\u00edmico6) Solu\u00e7\u00f5es\", replace: \"1) Introdu\u00e7\u00e3o; 2) 1a Lei da Termodin\u00e2mica 3) 2a e 3a Leis da Termodin\u00e2mica 4) Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase gasosa; 5) Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase condensada;\" },\n  { find: \"1- 1a Lei da Termodin\u00e2mica: sistema e vizinhan\u00e7a; transfer\u00eancia de energia; energia de um sistema; energia como uma fun\u00e7\u00e3o de estado; trabalho; sistema fechado; propriedades extensivas e intensivas; sistema aberto; entalpia; estado estacion\u00e1rio; capacidade t\u00e9rmica a volume constante; capacidade t\u00e9rmica a volume constante; equa\u00e7\u00e3o de estado, gases n\u00e3o-ideais; expans\u00e3o e compress\u00e3o adiab\u00e1tica; entalpias de forma\u00e7\u00e3o; varia\u00e7\u00e3o de entalpia em rea\u00e7\u00f5es qu\u00edmicas; varia\u00e7\u00e3o de temperatura associadas \u00e0 rea\u00e7\u00f5es qu\u00edmicas em sistemas adiab\u00e1ticos.2- 2a e 3a Leis da Termodin\u00e2mica: Entropia como fun\u00e7\u00e3o de estado; varia\u00e7\u00f5es de entropia associadas \u00e0 varia\u00e7\u00f5es de temperatura e press\u00e3o; varia\u00e7\u00f5es de entropia em rea\u00e7\u00f5es qu\u00edmicas; terceiro princ\u00edpio da termodin\u00e2mica.3- Rela\u00e7\u00f5es entre Propriedades Termodin\u00e2micas: As fun\u00e7\u00f5es A e G; potencial qu\u00edmico; grandezas parciais molares; rela\u00e7\u00f5es entre propriedades derivadas de U, H, A e G; g\u00e1s ideal; entropia de mistura; capacidade t\u00e9rmica; varia\u00e7\u00e3o de capacidade t\u00e9rmica; Rela\u00e7\u00e3o P-T isoentr\u00f3pica; compress\u00e3o isoentr\u00f3pica de s\u00f3lidos.4- Equil\u00edbrio: Condi\u00e7\u00f5es de equil\u00edbrio; equil\u00edbrio de fases; varia\u00e7\u00e3o de press\u00e3o de equil\u00edbrio com a temperatura; equa\u00e7\u00e3o de Clapeyron; varia\u00e7\u00e3o da press\u00e3o de vapor de uma fase condensada com a press\u00e3o total aplicada; varia\u00e7\u00e3o da press\u00e3o de vapor com tamanho de part\u00edcula.5- Equil\u00edbrio Qu\u00edmico: atividade; equil\u00edbrio qu\u00edmico; equil\u00edbrio em fase gasosa; equil\u00edbrio s\u00f3lido-vapor; fontes de informa\u00e7\u00e3o em DGo; diagrama de Ellingham; varia\u00e7\u00e3o da constante de equil\u00edbrio com a temperatura; gases dissolvidos em metais (Lei de Sievert); equil\u00edbrio qu\u00edmico e temperatura adiab\u00e1tica de chama.6- Solu\u00e7\u00f5es: grandezas parciais molares relativas; entropia de mistura - solu\u00e7\u00e3o ideal; entalpia de mistura  solu\u00e7\u00e3o ideal; Solu\u00e7\u00f5es n\u00e3o-ideais; rela\u00e7\u00e3o de Gibbs-Duhem; solu\u00e7\u00f5es regulares.\", replace: \"1- Introdu\u00e7\u00e3o: sistema; vizinhan\u00e7as; fases; equil\u00edbrio; fronteiras adiab\u00e1ticas e diat\u00e9rmicas; processos revers\u00edveis e irrevers\u00edveis; estado termodin\u00e2mico; mudan\u00e7a de estado; processos c\u00edclicos; equa\u00e7\u00e3o de estado; calor; trabalho.2- A 1\u00aa lei de Termodin\u00e2mica: energia interna; capacidades t\u00e9rmicas; entalpia; entalpia de transforma\u00e7\u00e3o de fases; entalpia de forma\u00e7\u00e3o e de rea\u00e7\u00e3o; entalpia de rea\u00e7\u00e3o em fun\u00e7\u00e3o da temperatura (introdu\u00e7\u00e3o ao loop termodin\u00e2mico).3- A 2\u00aa e 3\u00aa leis da Termodin\u00e2mica: Dispers\u00e3o de energia e entropia; entropia no zero absoluto; entropia de rea\u00e7\u00e3o; entropia de rea\u00e7\u00e3o em fun\u00e7\u00e3o da temperatura; desigualdade de Clausius; crit\u00e9rios de espontaneidade e equil\u00edbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de rea\u00e7\u00e3o em fun\u00e7\u00e3o da temperatura; equa\u00e7\u00e3o de Gibbs-Helmholtz.4- Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um g\u00e1s ideal; press\u00e3o de equil\u00edbrio em sistemas metal-\u00f3xido-O2(g).5- Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase condensada: fugacidade; atividade termodin\u00e2mica; solu\u00e7\u00f5es e grandezas parciais molares; potencial qu\u00edmico; modelos de solu\u00e7\u00f5es; propriedades termodin\u00e2micas de excesso\" },\n  { find: \"O curso ser\u00e1 ministrado na forma de aulas expositivas.\", replace: \"Esta \u00e9 uma disciplina fundamental, exigindo dedica\u00e7\u00e3o individual para assimila\u00e7\u00e3o de defini\u00e7\u00f5es e conceitos. Isto envolve leitura concentrada e realiza\u00e7\u00e3o de exerc\u00edcios num\u00e9ricos.\" },\n  { find: \"Para os alunos que obtiverem 3,0\u2264NF<5,0, ser\u00e1 aplicada uma avalia\u00e7\u00e3o de recupera\u00e7\u00e3o (R), com pontua\u00e7\u00e3o de 0 a 10, que levar\u00e1 ao c\u00e1lculo da m\u00e9dia final (MF) atrav\u00e9s da seguinte express\u00e3o:MF=(NF+R)/2\", replace: \"Para a recupera\u00e7\u00e3o ser\u00e1 realizada uma prova escrita (PR) abrangendo toda a mat\u00e9ria lecionada no semestre, valendo de 0 (zero) a 10 (dez). M\u00e9dia final = (NF + PR)/2.\" },\n  { find: \"REFER\u00caNCIAS BIBLIOGR\u00c1FICAS1) Castellan, G., Fundamentos de F\u00edsico-Qu\u00edmica, vol.1, Livros T\u00e9cnicos e Cient\u00edficos Editora S.A., 1986.2) Atkins, P.W., Physical Chemistry, Oxford University Press, 1994.3) Pilla, L., F\u00edsico Qu\u00edmica, vol.1, Livros T\u00e9cnicos e Cient\u00edficos Editora S.A., 1979.4) Moore, W.J., F\u00edsico Qu\u00edmica, vol.1, Editora Edgard Bl\u00fccher Ltda, 1976.5) Darken, L. & Gurry, R., Physical Chemistry f Metals, McGraw-Hill Book Company Inc., 1953.6) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995.7) Ragone, D.V., Thermodynamics of Materials, vol.1, John Wiley & Sons Inc., 1995. (Livro Texto)8) Swalin, R.A. Thermodynamics of Solids, John Wiley & Sons, 1972.9) Shoemaker, D.P., Garland, C.W., Nibler, J.W., Experiments in Physical Chemistry, McGraw-Hill Book Company, 1989.\", replace: \"1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. F\u00edsico-Qu\u00edmica, Livros T\u00e9cnicos e Cient\u00edficos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.\" },\n  { find: \"LOQ4098 -  Fundamentos de Qu\u00edmica para Engenharia II (Requisito fraco)\", replace: \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito fraco)\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find.substring(0, 60));\n  }\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "# Apply each literal text replacement with Find/Replace so that\n# sibling runs (e.g. line breaks) are left untouched.\n$d = $word.ActiveDocument\n\n$replacements = New-Object System.Collections.ArrayList\n\n$find0 = @'\nAtiva\u00e7\u00e3o: 01/01/2011\n'@\n$repl0 = @'\nAtiva\u00e7\u00e3o: 01/01/2023\n'@\n[void]$replacements.Add(@{Find = $find0; Replace = $repl0})\n\n$find1 = @'\nCurso (semestre ideal): EM (4)\n'@\n$repl1 = @'\nCurso (semestre ideal): EF (4), EM (4)\n'@\n[void]$replacements.Add(@{Find = $find1; Replace = $repl1})\n\n$find2 = @'\nO objetivo desta disciplina \u00e9 apresentar ao aluno fundamentos de F\u00edsico-Qu\u00edmica, mais especificamente Termodin\u00e2mica Qu\u00edmica, com foco em problemas encontrados na Engenharia de Materiais. Bastante \u00eanfase \u00e9 dada nos assuntos relativos a Termoqu\u00edmica e Equil\u00edbrio Qu\u00edmica, abordando equil\u00edbrio com fases condensadas.\n'@\n$repl2 = @'\nEsta disciplina visa apresentar fundamentos de termodin\u00e2mica aplicada \u00e0 \u00e1rea de ci\u00eancia e engenharia de materiais. Especial \u00eanfase \u00e9 dada \u00e0 energia na forma e calor para aquecimento de sistemas termodin\u00e2micos; c\u00e1lculos de varia\u00e7\u00e3o de entalpia; entropia e energia de Gibbs de elementos e compostos em mudan\u00e7as de estado; c\u00e1lculos de varia\u00e7\u00e3o de entalpia; entropia e energia de Gibbs de rea\u00e7\u00e3o; aplica\u00e7\u00e3o da propriedade  energia de Gibbs para avalia\u00e7\u00e3o de transforma\u00e7\u00f5es espont\u00e2neas e em equil\u00edbrio; fundamentos de termodin\u00e2mica de solu\u00e7\u00f5es; c\u00e1lculos de condi\u00e7\u00f5es de equil\u00edbrio em sistemas heterog\u00eaneos. Apresenta-se tamb\u00e9m as principais diferen\u00e7as entre esta disciplina e a disciplina de Termodin\u00e2mica de M\u00e1quinas.\n'@\n[void]$replacements.Add(@{Find = $find2; Replace = $repl2})\n\n$find3 = @'\n1) 1a Lei da Termodin\u00e2mica2) 2a e 3a Leis da Termodin\u00e2mica3) Rela\u00e7\u00f5es entre Propriedades Termodin\u00e2micas4) Equil\u00edbrio5) Equil\u00edbrio Qu\u00edmico6) Solu\u00e7\u00f5es\n'@\n$repl3 = @'\n1) Introdu\u00e7\u00e3o; 2) 1a Lei da Termodin\u00e2mica 3) 2a e 3a Leis da Termodin\u00e2mica 4) Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase gasosa; 5) Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase condensada;\n'@\n[void]$replacements.Add(@{Find = $find3; Replace = $repl3})\n\n$find4 = @'\n1- 1a Lei da Termodin\u00e2mica: sistema e vizinhan\u00e7a; transfer\u00eancia de energia; energia de um sistema; energia como uma fun\u00e7\u00e3o de estado; trabalho; sistema fechado; propriedades extensivas e intensivas; sistema aberto; entalpia; estado estacion\u00e1rio; capacidade t\u00e9rmica a volume constante; capacidade t\u00e9rmica a volume constante; equa\u00e7\u00e3o de estado, gases n\u00e3o-ideais; expans\u00e3o e compress\u00e3o adiab\u00e1tica; entalpias de forma\u00e7\u00e3o; varia\u00e7\u00e3o de entalpia em rea\u00e7\u00f5es qu\u00edmicas; varia\u00e7\u00e3o de temperatura associadas \u00e0 rea\u00e7\u00f5es qu\u00edmicas em sistemas adiab\u00e1ticos.2- 2a e 3a Leis da Termodin\u00e2mica: Entropia como fun\u00e7\u00e3o de estado; varia\u00e7\u00f5es de entropia associadas \u00e0 varia\u00e7\u00f5es de temperatura e press\u00e3o; varia\u00e7\u00f5es de entropia em rea\u00e7\u00f5es qu\u00edmicas; terceiro princ\u00edpio da termodin\u00e2mica.3- Rela\u00e7\u00f5es entre Propriedades Termodin\u00e2micas: As fun\u00e7\u00f5es A e G; potencial qu\u00edmico; grandezas parciais molares; rela\u00e7\u00f5es entre propriedades derivadas de U, H, A e G; g\u00e1s ideal; entropia de mistura; capacidade t\u00e9rmica; varia\u00e7\u00e3o de capacidade t\u00e9rmica; Rela\u00e7\u00e3o P-T isoentr\u00f3pica; compress\u00e3o isoentr\u00f3pica de s\u00f3lidos.4- Equil\u00edbrio: Condi\u00e7\u00f5es de equil\u00edbrio; equil\u00edbrio de fases; varia\u00e7\u00e3o de press\u00e3o de equil\u00edbrio com a temperatura; equa\u00e7\u00e3o de Clapeyron; varia\u00e7\u00e3o da press\u00e3o de vapor de uma fase condensada com a press\u00e3o total aplicada; varia\u00e7\u00e3o da press\u00e3o de vapor com tamanho de part\u00edcula.5- Equil\u00edbrio Qu\u00edmico: atividade; equil\u00edbrio qu\u00edmico; equil\u00edbrio em fase gasosa; equil\u00edbrio s\u00f3lido-vapor; fontes de informa\u00e7\u00e3o em DGo; diagrama de Ellingham; varia\u00e7\u00e3o da constante de equil\u00edbrio com a temperatura; gases dissolvidos em metais (Lei de Sievert); equil\u00edbrio qu\u00edmico e temperatura adiab\u00e1tica de chama.6- Solu\u00e7\u00f5es: grandezas parciais molares relativas; entropia de mistura - solu\u00e7\u00e3o ideal; entalpia de mistura  solu\u00e7\u00e3o ideal; Solu\u00e7\u00f5es n\u00e3o-ideais; rela\u00e7\u00e3o de Gibbs-Duhem; solu\u00e7\u00f5es regulares.\n'@\n$repl4 = @'\n1- Introdu\u00e7\u00e3o: sistema; vizinhan\u00e7as; fases; equil\u00edbrio; fronteiras adiab\u00e1ticas e diat\u00e9rmicas; processos revers\u00edveis e irrevers\u00edveis; estado termodin\u00e2mico; mudan\u00e7a de estado; processos c\u00edclicos; equa\u00e7\u00e3o de estado; calor; trabalho.2- A 1\u00aa lei de Termodin\u00e2mica: energia interna; capacidades t\u00e9rmicas; entalpia; entalpia de transforma\u00e7\u00e3o de fases; entalpia de forma\u00e7\u00e3o e de rea\u00e7\u00e3o; entalpia de rea\u00e7\u00e3o em fun\u00e7\u00e3o da temperatura (introdu\u00e7\u00e3o ao loop termodin\u00e2mico).3- A 2\u00aa e 3\u00aa leis da Termodin\u00e2mica: Dispers\u00e3o de energia e entropia; entropia no zero absoluto; entropia de rea\u00e7\u00e3o; entropia de rea\u00e7\u00e3o em fun\u00e7\u00e3o da temperatura; desigualdade de Clausius; crit\u00e9rios de espontaneidade e equil\u00edbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de rea\u00e7\u00e3o em fun\u00e7\u00e3o da temperatura; equa\u00e7\u00e3o de Gibbs-Helmholtz.4- Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um g\u00e1s ideal; press\u00e3o de equil\u00edbrio em sistemas metal-\u00f3xido-O2(g).5- Equil\u00edbrio heterog\u00eaneo: composi\u00e7\u00e3o vari\u00e1vel da fase condensada: fugacidade; atividade termodin\u00e2mica; solu\u00e7\u00f5es e grandezas parciais molares; potencial qu\u00edmico; modelos de solu\u00e7\u00f5es; propriedades termodin\u00e2micas de excesso\n'@\n[void]$replacements.Add(@{Find = $find4; Replace = $repl4})\n\n$find5 = @'\nO curso ser\u00e1 ministrado na forma de aulas expositivas.\n'@\n$repl5 = @'\nEsta \u00e9 uma disciplina fundamental, exigindo dedica\u00e7\u00e3o individual para assimila\u00e7\u00e3o de defini\u00e7\u00f5es e conceitos. Isto envolve leitura concentrada e realiza\u00e7\u00e3o de exerc\u00edcios num\u00e9ricos.\n'@\n[void]$replacements.Add(@{Find = $find5; Replace = $repl5})\n\n$find6 = @'\nPara os alunos que obtiverem 3,0\u2264NF<5,0, ser\u00e1 aplicada uma avalia\u00e7\u00e3o de recupera\u00e7\u00e3o (R), com pontua\u00e7\u00e3o de 0 a 10, que levar\u00e1 ao c\u00e1lculo da m\u00e9dia final (MF) atrav\u00e9s da seguinte express\u00e3o:MF=(NF+R)/2\n'@\n$repl6 = @'\nPara a recupera\u00e7\u00e3o ser\u00e1 realizada uma prova escrita (PR) abrangendo toda a mat\u00e9ria lecionada no semestre, valendo de 0 (zero) a 10 (dez). M\u00e9dia final = (NF + PR)/2.\n'@\n[void]$replacements.Add(@{Find = $find6; Replace = $repl6})\n\n$find7 = @'\nREFER\u00caNCIAS BIBLIOGR\u00c1FICAS1) Castellan, G., Fundamentos de F\u00edsico-Qu\u00edmica, vol.1, Livros T\u00e9cnicos e Cient\u00edficos Editora S.A., 1986.2) Atkins, P.W., Physical Chemistry, Oxford University Press, 1994.3) Pilla, L., F\u00edsico Qu\u00edmica, vol.1, Livros T\u00e9cnicos e Cient\u00edficos Editora S.A., 1979.4) Moore, W.J., F\u00edsico Qu\u00edmica, vol.1, Editora Edgard Bl\u00fccher Ltda, 1976.5) Darken, L. & Gurry, R., Physical Chemistry f Metals, McGraw-Hill Book Company Inc., 1953.6) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995.7) Ragone, D.V., Thermodynamics of Materials, vol.1, John Wiley & Sons Inc., 1995. (Livro Texto)8) Swalin, R.A. Thermodynamics of Solids, John Wiley & Sons, 1972.9) Shoemaker, D.P., Garland, C.W., Nibler, J.W., Experiments in Physical Chemistry, McGraw-Hill Book Company, 1989.\n'@\n$repl7 = @'\n1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. F\u00edsico-Qu\u00edmica, Livros T\u00e9cnicos e Cient\u00edficos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.\n'@\n[void]$replacements.Add(@{Find = $find7; Replace = $repl7})\n\n$find8 = @'\nLOQ4098 -  Fundamentos de Qu\u00edmica para Engenharia II (Requisito fraco)\n'@\n$repl8 = @'\nLOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito fraco)\n'@\n[void]$replacements.Add(@{Find = $find8; Replace = $repl8})\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $ok = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $($r.Find.Substring(0, [Math]::Min(60, $r.Find.Length)))\"\n    }\n}\n\nWrite-Output \"done\""}
